# Generate Report for Handback
# Regenerating the handback status report updates the "Correspond Handoff
# Datetime" and "Correspond Handback DateTime" columns. The 6042bfb4-... and
# a3e6584e-... rows shared the same timestamp text, so every cell that held
# that timestamp moves forward to the newly generated value.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 02:17:41"
$wsZhCn.Range("E4").Value = "2016-03-21 02:17:41"
$wsZhCn.Range("H3").Value = "2016-03-21 02:18:01"
$wsZhCn.Range("H4").Value = "2016-03-21 02:18:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 02:17:44"
$wsDeDe.Range("E4").Value = "2016-03-21 02:17:44"
$wsDeDe.Range("H3").Value = "2016-03-21 02:18:08"
$wsDeDe.Range("H4").Value = "2016-03-21 02:18:08"
